# Updated cryptos list on Thu Jul 25 09:10:13 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns for each coin row with
# newly-scraped figures, and re-sorts rows 45/46 (RenderToken now ranks
# above Hedera).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column holds strings that look numeric ("564.95", "6.00", ...).
# A plain .Value assignment lets Excel auto-coerce those into real numbers
# (losing the original text formatting, e.g. trailing zeros). Force the
# column to Text first so every write below round-trips as a string, then
# restore the default "Normal" style afterwards so no stray formatting is
# left behind on the cells.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# --- Rows 45/46 swapped places (RenderToken moved above Hedera) ---------
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "6.00"
$ws.Range("E45").Value = "  -6.84%  "

$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "0.0660"
$ws.Range("E46").Value = "  -4.85%  "

# --- Price / Volume(1h) refresh for the remaining rows -------------------
$ws.Range("D2").Value = "64.383.48"
$ws.Range("E2").Value = "  -2.97%  "

$ws.Range("D3").Value = "3.178.63"
$ws.Range("E3").Value = "  -8.13%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "564.95"
$ws.Range("E5").Value = "  -3.88%  "

$ws.Range("D6").Value = "170.68"
$ws.Range("E6").Value = "  -3.10%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -0.54%  "

$ws.Range("D9").Value = "3.173.67"
$ws.Range("E9").Value = "  -8.28%  "

$ws.Range("D10").Value = "0.125"
$ws.Range("E10").Value = "  -6.36%  "

$ws.Range("D11").Value = "6.64"
$ws.Range("E11").Value = "  -4.65%  "

$ws.Range("D12").Value = "0.397"
$ws.Range("E12").Value = "  -4.86%  "

$ws.Range("D13").Value = "3.725.00"
$ws.Range("E13").Value = "  -8.28%  "

$ws.Range("D14").Value = "0.136"
$ws.Range("E14").Value = "  +0.97%  "

$ws.Range("D15").Value = "27.44"
$ws.Range("E15").Value = "  -8.72%  "

$ws.Range("D16").Value = "64.331.49"
$ws.Range("E16").Value = "  -2.95%  "

$ws.Range("E17").Value = "  -5.20%  "

$ws.Range("D18").Value = "3.176.36"
$ws.Range("E18").Value = "  -8.20%  "

$ws.Range("D19").Value = "5.75"
$ws.Range("E19").Value = "  -3.63%  "

$ws.Range("E20").Value = "  -5.68%  "

$ws.Range("D21").Value = "353.31"
$ws.Range("E21").Value = "  -5.53%  "

$ws.Range("D22").Value = "7.20"
$ws.Range("E22").Value = "  -5.30%  "

$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("D24").Value = "69.09"
$ws.Range("E24").Value = "  -5.90%  "

$ws.Range("E25").Value = "  -5.53%  "

$ws.Range("D26").Value = "0.505"
$ws.Range("E26").Value = "  -5.79%  "

$ws.Range("D27").Value = "9.55"
$ws.Range("E27").Value = "  -3.74%  "

$ws.Range("E28").Value = "  -1.18%  "

$ws.Range("D29").Value = "1.02"
$ws.Range("E29").Value = "  +2.22%  "

$ws.Range("D30").Value = "5.62"
$ws.Range("E30").Value = "  -4.00%  "

$ws.Range("E32").Value = "  -5.13%  "

$ws.Range("D33").Value = "22.11"
$ws.Range("E33").Value = "  -6.79%  "

$ws.Range("E34").Value = "  -4.79%  "

$ws.Range("D35").Value = "6.63"
$ws.Range("E35").Value = "  -5.56%  "

$ws.Range("E36").Value = "  -7.21%  "

$ws.Range("D37").Value = "155.42"
$ws.Range("E37").Value = "  -3.66%  "

$ws.Range("D38").Value = "0.813"
$ws.Range("E38").Value = "  -7.95%  "

$ws.Range("D39").Value = "25.88"
$ws.Range("E39").Value = "  -9.11%  "

$ws.Range("E40").Value = "  -2.52%  "

$ws.Range("D41").Value = "1.71"
$ws.Range("E41").Value = "  -5.52%  "

$ws.Range("D42").Value = "2.612.43"
$ws.Range("E42").Value = "  -5.47%  "

$ws.Range("D43").Value = "4.19"
$ws.Range("E43").Value = "  -6.86%  "

$ws.Range("D44").Value = "39.68"
$ws.Range("E44").Value = "  -0.88%  "

$ws.Range("D47").Value = "24.06"
$ws.Range("E47").Value = "  -5.20%  "

$ws.Range("D48").Value = "325.17"
$ws.Range("E48").Value = "  -3.90%  "

$ws.Range("E49").Value = "  -7.28%  "

$ws.Range("E50").Value = "  -1.07%  "

$ws.Range("E51").Value = "  -0.10%  "

# Restore the default cell style now that every Price cell holds its text.
$priceRange.Style = "Normal"
